# edit.ps1
# Applies the diff: adds <w:proofErr> spellcheck markers around a few
# foreign/compound words, splits "Entry Condition:" / "Exit Condition:"
# into multiple runs, inserts a new "Il listino dell'azienda contiene
# almeno un prodotto" paragraph in the Entry Condition cell, and replaces
# "catalogo" with a separately-run "listino" in the Exit Condition cell.

$d = $word.ActiveDocument

function Wrap-Fragment([string]$innerXml) {
    $pre = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
    $post = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pre + $innerXml + $post
}

function Replace-ParagraphContainingText([string]$searchText, [string]$newParagraphXml) {
    $rng = $d.Content.Duplicate
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    $para = $rng.Paragraphs(1).Range
    $para.InsertXML((Wrap-Fragment $newParagraphXml))
}

# 1) "RFR2 - ModificaProdotto" title -> wrap "ModificaProdotto" run with proofErr
Replace-ParagraphContainingText "ModificaProdotto" @'
<w:p w14:paraId="7E5F1043" w14:textId="3E0F6D71" w:rsidR="00DB31FD" w:rsidRPr="009B7456" w:rsidRDefault="003860BE" w:rsidP="004F1535"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">RFR2 - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001712DA"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ModificaProdotto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

# 2) "Entry Condition:" -> split into "Entry " / "Condition" (proofErr) / ":"
Replace-ParagraphContainingText "Entry Condition:" @'
<w:p w14:paraId="0041C0E7" w14:textId="7D1D9858" w:rsidR="00DB31FD" w:rsidRPr="00DB31FD" w:rsidRDefault="00DB31FD"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00DB31FD"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Entry </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Condition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r></w:p>
'@

# 3) Entry-condition cell: keep " dei prodotti" paragraph, add new paragraph
#    "Il listino dell'azienda contiene almeno un prodotto" right after it.
Replace-ParagraphContainingText "dei prodotti" @'
<w:p w14:paraId="75E6E1BE" w14:textId="16B19B68" w:rsidR="00464621" w:rsidRPr="00464621" w:rsidRDefault="00464621"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">L’Azienda visualizza il </w:t></w:r><w:r w:rsidR="00F11AC5"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>listino</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> dei prodotti</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Il listino dell’azienda contiene almeno un prodotto</w:t></w:r></w:p>
'@

# 4) "L'Azienda compila il form, inserendo..." -> split "compila il form" so
#    "form" gets proofErr spellStart/spellEnd markers.
Replace-ParagraphContainingText "compila il form" @'
<w:p w14:paraId="1D062BF4" w14:textId="4DBBF2C2" w:rsidR="009B7456" w:rsidRPr="00464621" w:rsidRDefault="009B7456" w:rsidP="00464621"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr><w:ind w:left="373"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">L’Azienda </w:t></w:r><w:r w:rsidR="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">compila il </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>form</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> inserendo il nuovo</w:t></w:r><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">nome, </w:t></w:r><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>prezzo, foto e descrizione.</w:t></w:r></w:p>
'@

# 5) "Il sistema visualizza il form per la modifica del prodotto" -> split
#    out "form" with proofErr markers.
Replace-ParagraphContainingText "il form per la modifica del prodotto" @'
<w:p w14:paraId="6A21D76E" w14:textId="1582DA7D" w:rsidR="0013656A" w:rsidRPr="00464621" w:rsidRDefault="0013656A" w:rsidP="00464621"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:ind w:left="369"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">Il </w:t></w:r><w:r w:rsidR="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">istema visualizza </w:t></w:r><w:r w:rsidR="00464621" w:rsidRPr="00464621"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">il </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>form</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> per la modifica del prodotto</w:t></w:r></w:p>
'@

# 6) "Exit Condition:" -> split into "Exit " / "Condition" (proofErr) / ":"
Replace-ParagraphContainingText "Exit Condition:" @'
<w:p w14:paraId="483AA62F" w14:textId="2481F49B" w:rsidR="00DB31FD" w:rsidRPr="00DB31FD" w:rsidRDefault="00DB31FD" w:rsidP="00DB31FD"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00DB31FD"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Exit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Condition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r></w:p>
'@

# 7) "...ha modificato il prodotto presente nel catalogo..." -> split off
#    trailing word and replace it with a separately-run "listino".
Replace-ParagraphContainingText "ha modificato il prodotto presente nel catalogo" @'
<w:p w14:paraId="0C441CDA" w14:textId="2E729EA6" w:rsidR="00DB31FD" w:rsidRPr="00DB31FD" w:rsidRDefault="009B7456"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="003F2C89"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">L’Azienda </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">ha modificato il prodotto presente nel </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>listino</w:t></w:r><w:r w:rsidR="00F11AC5"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> e visualizza i dettagli del prodotto</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

Write-Host "All edits applied"
